$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("2023-12-08 11:01:56", 0.0006000000000000001),
    @("2023-12-08 11:02:39", 0.002),
    @("2023-12-08 11:03:16", 0.0026),
    @("2023-12-08 11:03:20", 0.0004)
)

$startRow = 100
for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $startRow + $i
    $ws.Cells.Item($row, 1).Value = $data[$i][0]
    $ws.Cells.Item($row, 2).Value = $data[$i][1]
}
